# Updates the "cryptos" price/volume table with refreshed market data.
# Cells whose new text would otherwise be auto-parsed by Excel as a plain
# number (losing formatting like trailing/leading zeros, e.g. "56.90" ->
# 56.9) are written with a leading apostrophe so Excel stores them as text,
# matching the original inline-string (text) representation of column D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.932.14"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "3.509.71"
$ws.Range("E3").Value = "  -4.61%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'578.91"
$ws.Range("D6").Value = "'175.21"
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "3.500.75"
$ws.Range("E8").Value = "  -4.65%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -7.33%  "
$ws.Range("D11").Value = "'6.59"
$ws.Range("E11").Value = "  +5.57%  "
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "'47.28"
$ws.Range("E13").Value = "  -5.36%  "
$ws.Range("D14").Value = "'0.0000278"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").Value = "'675.66"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "4.068.99"
$ws.Range("E17").Value = "  -4.67%  "
$ws.Range("D18").Value = "3.511.37"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").Value = "68.832.89"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D21").Value = "'17.58"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("D22").Value = "'11.17"
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("D23").Value = "'0.906"
$ws.Range("E23").Value = "  -4.24%  "
$ws.Range("E24").Value = "  -8.75%  "
$ws.Range("D25").Value = "'98.37"
$ws.Range("E25").Value = "  -5.20%  "
$ws.Range("E26").Value = "  -4.31%  "
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  -6.65%  "
$ws.Range("D30").Value = "'9.44"
$ws.Range("E30").Value = "  -8.37%  "
$ws.Range("D31").Value = "'33.03"
$ws.Range("E31").Value = "  -7.05%  "
$ws.Range("D32").Value = "'8.76"
$ws.Range("E32").Value = "  -5.06%  "
$ws.Range("E33").Value = "  -7.71%  "
$ws.Range("D34").Value = "'7.33"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  -6.34%  "
$ws.Range("D36").Value = "'573.23"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").Value = "'3.62"
$ws.Range("E37").Value = "  -14.37%  "
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("D40").Value = "'56.90"
$ws.Range("E40").Value = "  -5.58%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.137"
$ws.Range("E42").Value = "  -5.08%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0440"
$ws.Range("E43").Value = "  -4.88%  "
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("D45").Value = "3.418.50"
$ws.Range("E45").Value = "  -8.61%  "
$ws.Range("D46").Value = "'33.41"
$ws.Range("E46").Value = "  -6.46%  "
$ws.Range("D47").Value = "0.0₃0705"
$ws.Range("E47").Value = "  -9.09%  "
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("E49").Value = "  -7.55%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "'133.35"
$ws.Range("E51").Value = "  -0.66%  "
